$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BlockTypes")

# Step 1: set BB (image name) and BD (index) for each new row first
$ws.Range("BB16").Value = "blocksheet1_01.png"
$ws.Range("BD16").Value = 0
$ws.Range("BB17").Value = "blocksheet1_02.png"
$ws.Range("BD17").Formula = "=BD16+1"
$ws.Range("BB18").Value = "blocksheet1_03.png"
$ws.Range("BD18").Formula = "=BD17+1"
$ws.Range("BB20").Value = "blocksheet1_04.png"
$ws.Range("BD20").Value = 3
$ws.Range("BB21").Value = "blocksheet1_05.png"
$ws.Range("BD21").Value = 4
$ws.Range("BB23").Value = "blocksheet1_06.png"
$ws.Range("BD23").Value = 5
$ws.Range("BB24").Value = "blocksheet1_07.png"
$ws.Range("BD24").Value = 6
$ws.Range("BB27").Value = "blocksheet1_09.png"
$ws.Range("BD27").Value = 7
$ws.Range("BB29").Value = "blocksheet1_13.png"
$ws.Range("BD29").Value = 8
$ws.Range("BB30").Value = "blocksheet1_14.png"
$ws.Range("BD30").Value = 9
$ws.Range("BB32").Value = "blocksheet1_10.png"
$ws.Range("BD32").Value = 10
$ws.Range("BB33").Value = "blocksheet1_11.png"
$ws.Range("BD33").Value = 11
$ws.Range("BB34").Value = "blocksheet1_16.png"
$ws.Range("BD34").Value = 12
$ws.Range("BB35").Value = "blocksheet1_15.png"
$ws.Range("BD35").Value = 13
$ws.Range("BB36").Value = "blocksheet1_17.png"
$ws.Range("BD36").Value = 14
$ws.Range("BB46").Value = "blocksheet0-_08.png"
$ws.Range("BD46").Value = 7

# Step 2: set BC (column index) last, for all rows, so dependent formulas recalc correctly
$ws.Range("BC16").Value = 21
$ws.Range("BC17").Value = 21
$ws.Range("BC18").Value = 21
$ws.Range("BC20").Value = 21
$ws.Range("BC21").Value = 21
$ws.Range("BC23").Value = 21
$ws.Range("BC24").Value = 21
$ws.Range("BC27").Value = 21
$ws.Range("BC29").Value = 21
$ws.Range("BC30").Value = 21
$ws.Range("BC32").Value = 21
$ws.Range("BC33").Value = 21
$ws.Range("BC34").Value = 21
$ws.Range("BC35").Value = 21
$ws.Range("BC36").Value = 21
$ws.Range("BC46").Value = 20

# Step 3: AH5 formula was edited by hand in Excel (breaking the shared formula),
# producing an explicit formula text identical to the AH4 master but anchored to row 5
# (quirk: the original author referenced BB16 instead of BB5, preserved verbatim).
$ws.Range("AH5").Formula = "=IF(ISBLANK(E5),`"`",`"<>%  ===== `"&E5&`" ===========<>`"&AE5&`").<>inputs(`"&E5&`", `"&Z5&`").<>outputs(`"&E5&`", `"&AD5&`").<>input_types(`"&E5&`", `"&AP5&`").<>output_types(`"&E5&`", `"&AQ5&`").<>`"&IF(OR(ISBLANK(E5),ISBLANK(BB16)),`"`",`"image_name(`"&E5&`", '/img/blocks/`"&BB16&`"').<>`"))"

Write-Output "done"
